# Generate Report for Handback
#
# - Overview sheet: status text "Ready for handoff" -> "Handed back: in sync with en-US"
#   for both the zh-cn and de-de status columns (B/C) on both data rows.
# - zh-cn / de-de detail sheets: same status text change in column C; new
#   "Latest Target File" / "Latest Handback File" columns (F/G) populated
#   (re-using the already-known source/target file names & links), and the
#   "Latest Handback DateTime" (H) filled in for de-de (finished syncing)
#   while zh-cn keeps its placeholder (still pending).

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---- Overview sheet ------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("B2").Value = $handedBack
$overview.Range("C2").Value = $handedBack
$overview.Range("B3").Value = $handedBack
$overview.Range("C3").Value = $handedBack

# ---- zh-cn sheet -----------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C2").Value = $handedBack
$zhcn.Range("C3").Value = $handedBack
$zhcn.Range("H2").Value = "2016-03-12 14:48:05"
$zhcn.Range("H3").Value = "2016-03-12 14:48:05"

$zhcn.Hyperlinks.Add($zhcn.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/f9418a34f8018175bf1c07b7b0aaddad1c6567b1/e2e/8c50de45-616d-4b0f-9a5b-ab47a1647522.md", "", "", "8c50de45-616d-4b0f-9a5b-ab47a1647522.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82583653b107bddaafb8eaa3883b212fc92ef6a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/8c50de45-616d-4b0f-9a5b-ab47a1647522.61a71177ba8f505012c78f20b07b71b5dd019a8e.zh-cn.xlf", "", "", "8c50de45-616d-4b0f-9a5b-ab47a1647522.61a71177ba8f505012c78f20b07b71b5dd019a8e.zh-cn.xlf")
$zhcn.Hyperlinks.Add($zhcn.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/f9418a34f8018175bf1c07b7b0aaddad1c6567b1/e2e/f28a17c2-0989-40c6-852a-541543414afe.md", "", "", "f28a17c2-0989-40c6-852a-541543414afe.md")
$zhcn.Hyperlinks.Add($zhcn.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/82583653b107bddaafb8eaa3883b212fc92ef6a1/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/f28a17c2-0989-40c6-852a-541543414afe.99985fb90f70a4870b2223f4a87f5a18236a1af3.zh-cn.xlf", "", "", "f28a17c2-0989-40c6-852a-541543414afe.99985fb90f70a4870b2223f4a87f5a18236a1af3.zh-cn.xlf")

# ---- de-de sheet -----------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C2").Value = $handedBack
$dede.Range("C3").Value = $handedBack
$dede.Range("H2").Value = "2016-03-12 14:48:11"
$dede.Range("H3").Value = "2016-03-12 14:48:11"

$dede.Hyperlinks.Add($dede.Range("F2"), "https://github.com/OpenLocalizationTest/oltest/blob/f9418a34f8018175bf1c07b7b0aaddad1c6567b1/e2e/8c50de45-616d-4b0f-9a5b-ab47a1647522.md", "", "", "8c50de45-616d-4b0f-9a5b-ab47a1647522.md")
$dede.Hyperlinks.Add($dede.Range("G2"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c22268b03b1d26efbb1475eecf9966e22f663320/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/8c50de45-616d-4b0f-9a5b-ab47a1647522.61a71177ba8f505012c78f20b07b71b5dd019a8e.de-de.xlf", "", "", "8c50de45-616d-4b0f-9a5b-ab47a1647522.61a71177ba8f505012c78f20b07b71b5dd019a8e.de-de.xlf")
$dede.Hyperlinks.Add($dede.Range("F3"), "https://github.com/OpenLocalizationTest/oltest/blob/f9418a34f8018175bf1c07b7b0aaddad1c6567b1/e2e/f28a17c2-0989-40c6-852a-541543414afe.md", "", "", "f28a17c2-0989-40c6-852a-541543414afe.md")
$dede.Hyperlinks.Add($dede.Range("G3"), "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/c22268b03b1d26efbb1475eecf9966e22f663320/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/f28a17c2-0989-40c6-852a-541543414afe.99985fb90f70a4870b2223f4a87f5a18236a1af3.de-de.xlf", "", "", "f28a17c2-0989-40c6-852a-541543414afe.99985fb90f70a4870b2223f4a87f5a18236a1af3.de-de.xlf")
